$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 3915.4443
$ws.Range("J17").Value = 3915.4443
$ws.Range("L17").Value = 11746.3329
$ws.Range("N17").Value = -12082.3329

# Row 18
$ws.Range("H18").Value = 2284.111
$ws.Range("I18").Value = 2284.111
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 2284.111
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -2000.111

# Row 32
$ws.Range("H32").Value = 15943.5
$ws.Range("I32").Value = 15943.5
$ws.Range("K32").Value = 15943.5
$ws.Range("M32").Value = -15617.5

# Row 40
$ws.Range("H40").Value = 3100.4
$ws.Range("I40").Value = 2499
$ws.Range("K40").Value = 2499
$ws.Range("M40").Value = -2324

# Row 70
$ws.Range("H70").Value = 933291.9399999999
$ws.Range("J70").Value = 2117009.2
$ws.Range("L70").Value = 6351027.600000001
$ws.Range("N70").Value = -6351567.600000001

# Row 73
$ws.Range("H73").Value = 933291.9399999999
$ws.Range("J73").Value = 2117009.2
$ws.Range("L73").Value = 6351027.600000001
$ws.Range("N73").Value = -6352899.600000001

# Row 88
$ws.Range("H88").Value = 2131.8125
$ws.Range("I88").Value = 1032.2858
$ws.Range("J88").Value = 2987
$ws.Range("K88").Value = 1032.2858
$ws.Range("L88").Value = 2987
$ws.Range("M88").Value = -626.2858000000001
$ws.Range("N88").Value = -3799

# Row 91
$ws.Range("H91").Value = 2131.8125
$ws.Range("I91").Value = 1032.2858
$ws.Range("J91").Value = 2987
$ws.Range("K91").Value = 1032.2858
$ws.Range("L91").Value = 2987
$ws.Range("M91").Value = 371.7141999999999
$ws.Range("N91").Value = -5795

# Row 132
$ws.Range("H132").Value = 1136.5161
$ws.Range("I132").Value = 916.0741
$ws.Range("K132").Value = 2748.2223
$ws.Range("M132").Value = -218.2223000000004

# Row 137
$ws.Range("H137").Value = 2927.625
$ws.Range("I137").Value = 2585.4666
$ws.Range("J137").Value = 3497.889
$ws.Range("K137").Value = 7756.399800000001
$ws.Range("L137").Value = 10493.667
$ws.Range("M137").Value = -5206.399800000001
$ws.Range("N137").Value = -15593.667

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 393.2
$ws.Range("I5").Value = 366.5
$ws.Range("K5").Value = 366.5
$ws.Range("M5").Value = -254.5

# Row 32
$ws.Range("H32").Value = 33807.438
$ws.Range("I32").Value = 32819.223
$ws.Range("K32").Value = 32819.223
$ws.Range("M32").Value = -32532.223

# Row 38
$ws.Range("H38").Value = 27450
$ws.Range("I38").Value = 15000
$ws.Range("J38").Value = 39900
$ws.Range("K38").Value = 15000
$ws.Range("L38").Value = 39900
$ws.Range("M38").Value = -14533
$ws.Range("N38").Value = -40834

# Row 74
$ws.Range("H74").Value = 2164.6875
$ws.Range("I74").Value = 782.0952
$ws.Range("K74").Value = 782.0952
$ws.Range("M74").Value = 91.90480000000002

# Row 77
$ws.Range("H77").Value = 2164.6875
$ws.Range("I77").Value = 782.0952
$ws.Range("K77").Value = 3910.476
$ws.Range("M77").Value = 457.5240000000003

# Row 88
$ws.Range("H88").Value = 1875.125
$ws.Range("J88").Value = 1857.4286
$ws.Range("L88").Value = 1857.4286
$ws.Range("N88").Value = -2669.4286

# Row 91
$ws.Range("H91").Value = 1875.125
$ws.Range("J91").Value = 1857.4286
$ws.Range("L91").Value = 1857.4286
$ws.Range("N91").Value = -4665.4286

# Row 110
$ws.Range("H110").Value = 7577152.5
$ws.Range("I110").Value = 10000842
$ws.Range("J110").Value = 3122.5
$ws.Range("K110").Value = 10000842
$ws.Range("L110").Value = 3122.5
$ws.Range("M110").Value = -9998797
$ws.Range("N110").Value = -7212.5

# Row 132
$ws.Range("H132").Value = 8249.083000000001
$ws.Range("I132").Value = 6339.077
$ws.Range("K132").Value = 19017.231
$ws.Range("M132").Value = -16487.231

# Row 140
$ws.Range("H140").Value = 50213.5
$ws.Range("J140").Value = 50213.5
$ws.Range("L140").Value = 50213.5
$ws.Range("N140").Value = -60573.5

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 393.2
$ws.Range("I4").Value = 366.5
$ws.Range("K4").Value = 366.5
$ws.Range("M4").Value = -251.5

# Row 94
$ws.Range("H94").Value = 1380.7241
$ws.Range("I94").Value = 1274.6666
$ws.Range("K94").Value = 1274.6666
$ws.Range("M94").Value = -823.6666

# Row 99
$ws.Range("H99").Value = 2800
$ws.Range("I99").Value = 3080
$ws.Range("J99").Value = 2100
$ws.Range("K99").Value = 3080
$ws.Range("L99").Value = 2100
$ws.Range("M99").Value = -1582
$ws.Range("N99").Value = -5096

# Row 105
$ws.Range("H105").Value = 250067500
$ws.Range("I105").Value = 250067500
$ws.Range("K105").Value = 250067500
$ws.Range("M105").Value = -250065753

$ws = $wb.Worksheets.Item("CRP")
# Row 19
$ws.Range("H19").Value = 1053.4286
$ws.Range("I19").Value = 421.25
$ws.Range("J19").Value = 1896.3334
$ws.Range("K19").Value = 421.25
$ws.Range("L19").Value = 1896.3334
$ws.Range("M19").Value = -251.25
$ws.Range("N19").Value = -2236.3334

# Row 21
$ws.Range("H21").Value = 9999.5
$ws.Range("J21").Value = 9999.5
$ws.Range("L21").Value = 9999.5
$ws.Range("N21").Value = -10469.5

# Row 24
$ws.Range("H24").Value = 1053.4286
$ws.Range("I24").Value = 421.25
$ws.Range("J24").Value = 1896.3334
$ws.Range("K24").Value = 421.25
$ws.Range("L24").Value = 1896.3334
$ws.Range("M24").Value = -251.25
$ws.Range("N24").Value = -2236.3334

# Row 31
$ws.Range("H31").Value = 35718764
$ws.Range("J31").Value = 7653.5
$ws.Range("L31").Value = 7653.5
$ws.Range("N31").Value = -8243.5

# Row 34
$ws.Range("H34").Value = 35718764
$ws.Range("J34").Value = 7653.5
$ws.Range("L34").Value = 7653.5
$ws.Range("N34").Value = -8057.5

# Row 105
$ws.Range("H105").Value = 1949.6
$ws.Range("I105").Value = 1187
$ws.Range("K105").Value = 1187
$ws.Range("M105").Value = 560

# Row 132
$ws.Range("H132").Value = 18572.469
$ws.Range("I132").Value = 746.1667
$ws.Range("K132").Value = 2238.5001
$ws.Range("M132").Value = 291.4998999999998

# Row 134
$ws.Range("H134").Value = 5396.304
$ws.Range("I134").Value = 4555.316
$ws.Range("K134").Value = 13665.948
$ws.Range("M134").Value = -11130.948

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 30272.56
$ws.Range("I4").Value = 34541.33
$ws.Range("K4").Value = 103623.99
$ws.Range("M4").Value = -103511.99

$ws = $wb.Worksheets.Item("GSM")
# Row 31
$ws.Range("H31").Value = 3491.1667
$ws.Range("I31").Value = 1236.75
$ws.Range("K31").Value = 1236.75
$ws.Range("M31").Value = -944.75

# Row 37
$ws.Range("H37").Value = 3491.1667
$ws.Range("I37").Value = 1236.75
$ws.Range("K37").Value = 1236.75
$ws.Range("M37").Value = -959.75

# Row 97
$ws.Range("H97").Value = 1790
$ws.Range("I97").Value = 2212
$ws.Range("J97").Value = 1262.5
$ws.Range("K97").Value = 2212
$ws.Range("L97").Value = 1262.5
$ws.Range("M97").Value = -1716
$ws.Range("N97").Value = -2254.5

# Row 113
$ws.Range("H113").Value = 2494.1428
$ws.Range("I113").Value = 2393.1667
$ws.Range("K113").Value = 2393.1667
$ws.Range("M113").Value = -223.1667000000002

# Row 122
$ws.Range("H122").Value = 1432.75
$ws.Range("J122").Value = 1212.9
$ws.Range("L122").Value = 3638.7
$ws.Range("N122").Value = -8538.700000000001

# Row 138
$ws.Range("H138").Value = 74000

$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 5272.07
$ws.Range("I136").Value = 4459.0938
$ws.Range("J136").Value = 7637.091
$ws.Range("K136").Value = 13377.2814
$ws.Range("L136").Value = 22911.273
$ws.Range("M136").Value = -10827.2814
$ws.Range("N136").Value = -28011.273

$ws = $wb.Worksheets.Item("WVR")
# Row 43
$ws.Range("H43").Value = 76666
$ws.Range("J43").Value = 64999
$ws.Range("L43").Value = 64999
$ws.Range("N43").Value = -65297

# Row 70
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

# Row 73
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

# Row 109
$ws.Range("H109").Value = 40937
$ws.Range("I109").Value = 18000
$ws.Range("K109").Value = 18000
$ws.Range("M109").Value = -16613

# Row 136
$ws.Range("H136").Value = 4225.85
$ws.Range("I136").Value = 1714.3334
$ws.Range("J136").Value = 7993.125
$ws.Range("K136").Value = 5143.0002
$ws.Range("L136").Value = 23979.375
$ws.Range("M136").Value = -2593.0002
$ws.Range("N136").Value = -29079.375
